$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.295.73'
$ws.Range('E2').Value = '  -2.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.383.20'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.03'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.44'
$ws.Range('E6').Value = '  +7.05%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.385.23'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.58'
$ws.Range('E10').Value = '  +2.60%  '
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.960.09'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.384.49'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.18'
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.406.86'
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.06'
$ws.Range('E19').Value = '  +6.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.80'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.40'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '375.88'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.568'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.518.19'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.64'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  +9.10%  '
$ws.Range('E28').Value = '  +21.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.75'
$ws.Range('E29').Value = '  +10.96%  '
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.14'
$ws.Range('E31').Value = '  +3.81%  '
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.416.67'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.41'
$ws.Range('E36').Value = '  +2.99%  '
$ws.Range('E37').Value = '  +5.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.97'
$ws.Range('E38').Value = '  +3.77%  '
$ws.Range('E39').Value = '  +5.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.45'
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0786'
$ws.Range('E41').Value = '  +3.98%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.22'
$ws.Range('E43').Value = '  +12.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.42'
$ws.Range('E44').Value = '  +3.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.46'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.761'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('E47').Value = '  +2.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.53'
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.94'
$ws.Range('E49').Value = '  +4.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.03'
$ws.Range('E50').Value = '  +12.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.900'
$ws.Range('E51').Value = '  +5.49%  '
